# Update statbar xpaths & diagnosis testcases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# column B -> C and column C -> D, leaving their values/widths intact.
$ws.Columns("B").Insert()

# New column B header + long-form stat query text.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Lymphoma']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Match the wrap-text style used by the other long-text cell in row 2 (A2).
$ws.Range("B2").WrapText = $true

# New column B should share column A's width (Excel normally carries the
# left-neighbour's width onto an inserted column); columns C/D (former
# B/C) keep their original widths untouched.
$ws.Columns("B").ColumnWidth = 75

# Selection moves to the single new cell B2.
$ws.Range("B2").Select()
